$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 32-37 (Manufacturing .. Construction) were finished by Farid but the
# "BY" / "STATUS" columns were left blank - fill them in like the rows above.
$doneRows = 32..37
foreach ($r in $doneRows) {
    $ws.Cells.Item($r, 4).Value = "Farid"
    $ws.Cells.Item($r, 5).Value = "DONE"
}

# Row 46 (Industrial) is a brand new project that was missing - Fahri picked
# it up and finished it.
$ws.Cells.Item(46, 4).Value = "Fahri"
$ws.Cells.Item(46, 5).Value = "DONE"

# Move the viewport / selection to reflect where work left off.
[void]$ws.Range("F46").Select()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
